$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '28.176.07'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -1.37%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.804.13'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  -0.15%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '316.90'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.99%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.13%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5346'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.23%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3773'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -0.88%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.07476'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.60%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '42.06'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -0.93%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '1.096'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -1.82%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.21%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '6.215'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('E14').Value = '  -2.69%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '7.394'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.39%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '1.801.97'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.69%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '89.87'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.38%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.00001063'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.18%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.06506'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.96%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '17.40'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.91%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '1.0000'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.13%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.928'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.12%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '28.200.77'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -1.31%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '11.19'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.12%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.087'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.53%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '156.14'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -3.17%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '20.52'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.08%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.015.50'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +0.94%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.329'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -2.18%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '121.98'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.96%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.129'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +0.60%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.1101'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +7.78%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '5.599'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -1.88%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '3.624'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.89%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.07094'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +8.26%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.2228'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -3.36%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.02301'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.87%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '5.102'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +0.26%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '8.469'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.74%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.6176'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -2.23%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '11.13'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -2.97%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.183'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -2.49%  '
$ws.Range('B43').Value = 'WEMIXTOKEN'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.434'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +3.73%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '13.38'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -1.34%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '3.686'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.52%  '
$ws.Range('E46').Value = '  -2.58%  '
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.186'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.925'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -2.68%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.06820'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -1.52%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '71.76'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.04%  '
